$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "StatQuery" text (Cypher query) that replaces the previous generic
# stat query used across the CasesTab / SamplesTab / FilesTab rows.
$newStatQuery = @'
MATCH (p:program)<--(s:study)<-[*]-(c:case)<--(demo:demographic)
OPTIONAL MATCH (samp:sample)-->(c)
OPTIONAL MATCH (diag:diagnosis)-->(c)
OPTIONAL MATCH (f:file)-[*]->(c)
OPTIONAL MATCH (sf:file)-->(s)
WITH DISTINCT f, sf, samp AS samp, c, demo, diag, s, p
WHERE demo.breed IN ['Chinese Shar-Pei']
RETURN  
    count(distinct p) AS Programs,
    count(distinct s) AS Studies,
    count(distinct c) AS Cases,
    count(distinct samp) AS Samples,
    count(distinct f) AS `Case Files`,
    count(distinct sf) AS `Study Files`
'@

# The original query text ends without a trailing newline; strip the one
# introduced by the here-string syntax so the cell content matches exactly.
$newStatQuery = $newStatQuery.TrimEnd("`r", "`n")

# Column C ("StatQuery") for the three data rows (CasesTab, SamplesTab,
# FilesTab) all get updated to the new query text.
$ws.Range("C2").Value = $newStatQuery
$ws.Range("C3").Value = $newStatQuery
$ws.Range("C4").Value = $newStatQuery

# Update the saved view state: scrolled down so row 3 is at the top,
# zoomed to 85%, with B4 selected as the active cell.
$window = $excel.ActiveWindow
$window.Zoom = 85
$window.ScrollRow = 3
$window.ScrollColumn = 1
$null = $ws.Range("B4").Select()
